# AFDP-7308 Combine Transcribe and OCR processing into a single media processing module
#
# The "Transcribe Workflow Rules" sheet referenced a dedicated
# com.armedia.acm.services.transcribe.model.TranscribeBusinessProcessModel /
# TranscribeWorkflow business process. This edit repoints those two rule
# actions (rows 17-18, "Automatic Transcribe" / "Manual Transcribe") at the
# new, combined media-processing module instead:
#   - com.armedia.acm.services.transcribe.model.TranscribeBusinessProcessModel
#       -> com.armedia.acm.services.mediaengine.model.MediaEngineBusinessProcessModel
#   - $model: TranscribeBusinessProcessModel
#       -> $model: MediaEngineBusinessProcessModel
#   - TranscribeWorkflow -> MediaEngineWorkFlow (process started by both rules)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Business process model class used by the rule table's imports
$ws.Range("D3").Value = "com.armedia.acm.services.mediaengine.model.MediaEngineBusinessProcessModel"

# The workflow process name started by both the automatic and manual
# transcribe rules now points at the combined media engine workflow
$ws.Range("E17").Value = "MediaEngineWorkFlow"
$ws.Range("E18").Value = "MediaEngineWorkFlow"

# Global variable declaration line in the RuleTable header
$ws.Range("C14").Value = "`$model: MediaEngineBusinessProcessModel"

# The wrapped-text rows affected by the surrounding re-save reflow to new
# heights even though their own text is unchanged
$ws.Rows.Item(9).RowHeight = 158.4
$ws.Rows.Item(16).RowHeight = 115.2

# Restore the view further down the sheet, centred on the rule table that
# was just edited
$ws.Range("A12:E18").Select()

$wb.Saved = $false
